# Update odds values on Sheet1 (rows 6-14) to match the latest FlashScore
# refresh for 2025-02-18, per the commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = 3.7
$ws.Range("I6").Value = 2.2
$ws.Range("J6").Value = 4.3
$ws.Range("L6").Value = 2.82
$ws.Range("P6").Value = 2.4
$ws.Range("S6").Value = 2.5
$ws.Range("T6").Value = 1.47
$ws.Range("W6").Value = 4.4
$ws.Range("X6").Value = 1.17
$ws.Range("AA6").Value = 2.05
$ws.Range("AB6").Value = 1.7
$ws.Range("AC6").Value = 7.9
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 60
$ws.Range("AG6").Value = 45
$ws.Range("AH6").Value = 60
$ws.Range("AK6").Value = 17
$ws.Range("AL6").Value = 110
$ws.Range("AN6").Value = 5.8
$ws.Range("AP6").Value = 9
$ws.Range("AS6").Value = 40
$ws.Range("K7").Value = 1.8
$ws.Range("Q7").Value = 2.04
$ws.Range("R7").Value = 1.72
$ws.Range("AA7").Value = 2.25
$ws.Range("AB7").Value = 1.57
$ws.Range("AC7").Value = 5.5
$ws.Range("AI7").Value = 6
$ws.Range("AN7").Value = 7
$ws.Range("AS7").Value = 51
$ws.Range("G8").Value = 2.45
$ws.Range("I8").Value = 2.7
$ws.Range("L8").Value = 3.5
$ws.Range("AA8").Value = 1.8
$ws.Range("AB8").Value = 1.95
$ws.Range("AF8").Value = 23
$ws.Range("AN8").Value = 8.5
$ws.Range("AP8").Value = 11
$ws.Range("AQ8").Value = 29
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("S9").Value = 2.1
$ws.Range("T9").Value = 1.7
$ws.Range("W9").Value = 3.75
$ws.Range("X9").Value = 1.25
$ws.Range("G10").Value = 3.6
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 4.5
$ws.Range("L10").Value = 2.75
$ws.Range("AA10").Value = 1.95
$ws.Range("AB10").Value = 1.8
$ws.Range("AD10").Value = 17
$ws.Range("AL10").Value = 51
$ws.Range("AO10").Value = 9
$ws.Range("AQ10").Value = 17
$ws.Range("S11").Value = 2.15
$ws.Range("T11").Value = 1.67
$ws.Range("G13").Value = 1.4
$ws.Range("H13").Value = 4.33
$ws.Range("I13").Value = 7.5
$ws.Range("J13").Value = 1.87
$ws.Range("AD13").Value = 7
$ws.Range("AF13").Value = 9.5
$ws.Range("AN13").Value = 19
$ws.Range("AP13").Value = 21
$ws.Range("AA14").Value = 1.77
$ws.Range("AB14").Value = 1.87
